# "Adding comments to Example2..." -- Trello1.xlsx:
#   - rename the worksheet tab from "Sheet1" to "Trello Tasks"
#   - refresh the sample Trello card id stored in E2 with a new GUID
#   - leave the selection sitting on the cell that was just edited (E2)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab.
$ws.Name = "Trello Tasks"

# Update the Trello card id (was 89d78c8f-4bd2-4d8e-89e8-878173996911).
$ws.Range("E2").Value = "a6ecd523-39f9-47da-90e0-f9ff7b3691c4"

# Park the active selection on the cell we just changed.
$ws.Range("E2").Select()
